$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Learning_room").Name = "Example_1"
$wb.Worksheets.Item("Carpentry_workshop").Name = "Example_2"
$wb.Worksheets.Item("Office").Name = "Example_3"
$wb.Worksheets.Item("Emergeny_room").Name = "Example_4"
$wb.Worksheets.Item("Industrial_production_line").Name = "Example_5"

$wb.Worksheets.Item("Example_5").Activate()
